# Applies the edits described by the commit diff to the active document.
$d = $word.ActiveDocument

# 1) "- База данных российских магазинов и продаваемых в них цифровых товаров;"
#    (client bullet list) -> "- Кэш базы данных для работы в оффлайн-режиме;"
#    NOTE: must run before step 1b below, since that step introduces a second
#    occurrence of this same phrase (in the server bullet list) which would
#    otherwise also get matched here.
$d.Content.Find.Execute(
    "База данных российских магазинов и продаваемых в них цифровых товаров;",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Кэш базы данных для работы в оффлайн-режиме;",
    2) | Out-Null

# 1b) "- Готовая база данных товаров и магазинов;" (server bullet list)
#    -> "- База данных российских магазинов и продаваемых в них цифровых товаров;"
$d.Content.Find.Execute(
    "Готовая база данных товаров и магазинов",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "База данных российских магазинов и продаваемых в них цифровых товаров",
    2) | Out-Null

# 2) "- Система обслуживания и администрирования сервера;"
#    -> "- Система обслуживания и администрирования сервера (поддержка стабильной работы, uptime, внесение необходимых изменений в параметры сервера);"
$d.Content.Find.Execute(
    "Система обслуживания и администрирования сервера;",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Система обслуживания и администрирования сервера (поддержка стабильной работы, uptime, внесение необходимых изменений в параметры сервера);",
    2) | Out-Null

# 3) "- Программный интерфейс для управления сервером;"
#    -> "- Программный интерфейс для управления сервером (включение, отключение, перезагрузка, мониторинг нагрузки и активности);"
$d.Content.Find.Execute(
    "Программный интерфейс для управления сервером;",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Программный интерфейс для управления сервером (включение, отключение, перезагрузка, мониторинг нагрузки и активности);",
    2) | Out-Null

# 5) "Соисполнитель №1 ... программного и пользовательского интерфейса ..."
#    -> drop "пользовательского "
$d.Content.Find.Execute(
    "программного и пользовательского интерфейса",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "программного и интерфейса",
    2) | Out-Null

# 6) "Соисполнитель №2 - разработка системы поиска по программе, сбор необходимой
#     информации о товарах и интернет-магазинах, разработка системы обратной связи;"
#    -> "Соисполнитель №2 – разработка пользовательского интерфейса, разработка
#        системы поиска по программе, разработка системы обратной связи;"
$d.Content.Find.Execute(
    "Соисполнитель №2 - разработка системы поиска по программе, сбор необходимой информации о товарах и интернет-магазинах, разработка системы обратной связи;",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Соисполнитель №2 – разработка пользовательского интерфейса, разработка системы поиска по программе, разработка системы обратной связи;",
    2) | Out-Null
